$d = $word.ActiveDocument

# Locate the paragraph containing the "Requisitos" detail line, then remove
# the trailing blank paragraph plus the "Ver no Jupiter..." and
# "(c) 2020 ..." footer paragraphs that used to follow it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOM3212: Fen*menos de Transporte A (Requisito)*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
